# Change styles in reference.docx: give all headings a uniform color,
# shrink Heading 2 / Heading 3, and swap Heading 4 from bold to italic
# while un-italicizing Heading 5 (matching commit "Closes #5820").

$d = $word.ActiveDocument

# --- Heading 1: drop the themeShade="B5" darkening; use the plain
#     accent1 color (4F81BD) that the other headings already use. ---
$h1 = $d.Styles("Heading 1")
$h1.Font.Color = 12419407   ; # RGB(0x4F, 0x81, 0xBD) == accent1, no shade

# --- Heading 2: 16pt -> 14pt (both the Latin and complex-script size). ---
$h2 = $d.Styles("Heading 2")
$h2.Font.Size = 14
$h2.Font.SizeBi = 14

# --- Heading 3: 14pt -> 12pt (both the Latin and complex-script size). ---
$h3 = $d.Styles("Heading 3")
$h3.Font.Size = 12
$h3.Font.SizeBi = 12

# --- Heading 4: bold -> italic. ---
$h4 = $d.Styles("Heading 4")
$h4.Font.Bold = $false
$h4.Font.Italic = $true

# --- Heading 5: no longer italic. ---
$h5 = $d.Styles("Heading 5")
$h5.Font.Italic = $false
